$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation in provider/company names (comma -> period) ---
$ws.Range("E33").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F33").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E64").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("F64").Value = 'PITTER ROLANDO LJ. CERGNEUX MARIO M Y CERGNEUX DANIEL F  SH'
$ws.Range("E72").Value = 'FERNANDEZ MARIO H. GALLICET OSCAR M'
$ws.Range("E73").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("F73").Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range("E75").Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range("E98").Value = 'RICCOTTI. MARIANA EDITH'

# --- Reformat "Importe" amounts: "1.234,56" (es-AR) -> "1234.56" (plain) ---
# Values are stored as text; prefix with an apostrophe so Excel keeps them as text
# (otherwise Excel would coerce the plain-decimal form to a Number and drop trailing zeros).
$ws.Cells.Item(2, 8).Value = "'1020.00"
$ws.Cells.Item(3, 8).Value = "'384.00"
$ws.Cells.Item(4, 8).Value = "'1092.00"
$ws.Cells.Item(5, 8).Value = "'25000.50"
$ws.Cells.Item(6, 8).Value = "'90.00"
$ws.Cells.Item(7, 8).Value = "'164.00"
$ws.Cells.Item(8, 8).Value = "'1348.10"
$ws.Cells.Item(9, 8).Value = "'10260.00"
$ws.Cells.Item(10, 8).Value = "'930.95"
$ws.Cells.Item(11, 8).Value = "'19220.45"
$ws.Cells.Item(12, 8).Value = "'4620.07"
$ws.Cells.Item(13, 8).Value = "'180.00"
$ws.Cells.Item(14, 8).Value = "'49806.50"
$ws.Cells.Item(15, 8).Value = "'91545.20"
$ws.Cells.Item(16, 8).Value = "'27182.49"
$ws.Cells.Item(17, 8).Value = "'720.00"
$ws.Cells.Item(18, 8).Value = "'1387.65"
$ws.Cells.Item(19, 8).Value = "'5967.40"
$ws.Cells.Item(20, 8).Value = "'206.40"
$ws.Cells.Item(21, 8).Value = "'954.80"
$ws.Cells.Item(22, 8).Value = "'5392.24"
$ws.Cells.Item(23, 8).Value = "'450.00"
$ws.Cells.Item(24, 8).Value = "'438.00"
$ws.Cells.Item(25, 8).Value = "'890.00"
$ws.Cells.Item(26, 8).Value = "'1077.00"
$ws.Cells.Item(27, 8).Value = "'3300.00"
$ws.Cells.Item(28, 8).Value = "'53.09"
$ws.Cells.Item(29, 8).Value = "'285.00"
$ws.Cells.Item(30, 8).Value = "'620.00"
$ws.Cells.Item(31, 8).Value = "'1807.46"
$ws.Cells.Item(32, 8).Value = "'167.50"
$ws.Cells.Item(33, 8).Value = "'505.66"
$ws.Cells.Item(34, 8).Value = "'169.54"
$ws.Cells.Item(35, 8).Value = "'129.00"
$ws.Cells.Item(36, 8).Value = "'540.00"
$ws.Cells.Item(37, 8).Value = "'25973.93"
$ws.Cells.Item(38, 8).Value = "'319.20"
$ws.Cells.Item(39, 8).Value = "'333.50"
$ws.Cells.Item(40, 8).Value = "'191.00"
$ws.Cells.Item(41, 8).Value = "'3802.70"
$ws.Cells.Item(42, 8).Value = "'57822.13"
$ws.Cells.Item(43, 8).Value = "'33280.15"
$ws.Cells.Item(44, 8).Value = "'1235.83"
$ws.Cells.Item(45, 8).Value = "'1931.00"
$ws.Cells.Item(46, 8).Value = "'192.00"
$ws.Cells.Item(47, 8).Value = "'181.65"
$ws.Cells.Item(48, 8).Value = "'27709.08"
$ws.Cells.Item(49, 8).Value = "'54.00"
$ws.Cells.Item(50, 8).Value = "'2149.11"
$ws.Cells.Item(51, 8).Value = "'78.00"
$ws.Cells.Item(52, 8).Value = "'31.80"
$ws.Cells.Item(53, 8).Value = "'997.00"
$ws.Cells.Item(54, 8).Value = "'1290.50"
$ws.Cells.Item(55, 8).Value = "'25.30"
$ws.Cells.Item(56, 8).Value = "'306.00"
$ws.Cells.Item(57, 8).Value = "'435.50"
$ws.Cells.Item(58, 8).Value = "'1822.00"
$ws.Cells.Item(59, 8).Value = "'38.46"
$ws.Cells.Item(60, 8).Value = "'11520.00"
$ws.Cells.Item(61, 8).Value = "'4973.00"
$ws.Cells.Item(62, 8).Value = "'4184.00"
$ws.Cells.Item(63, 8).Value = "'300.08"
$ws.Cells.Item(64, 8).Value = "'405.00"
$ws.Cells.Item(65, 8).Value = "'1256.00"
$ws.Cells.Item(66, 8).Value = "'720.00"
$ws.Cells.Item(67, 8).Value = "'10345.80"
$ws.Cells.Item(68, 8).Value = "'7879.00"
$ws.Cells.Item(69, 8).Value = "'38.00"
$ws.Cells.Item(70, 8).Value = "'80.00"
$ws.Cells.Item(71, 8).Value = "'1340.00"
$ws.Cells.Item(72, 8).Value = "'160.00"
$ws.Cells.Item(73, 8).Value = "'2139.73"
$ws.Cells.Item(74, 8).Value = "'580.00"
$ws.Cells.Item(75, 8).Value = "'260.00"
$ws.Cells.Item(76, 8).Value = "'0.08"
$ws.Cells.Item(77, 8).Value = "'0.08"
$ws.Cells.Item(78, 8).Value = "'2535.00"
$ws.Cells.Item(79, 8).Value = "'4.80"
$ws.Cells.Item(80, 8).Value = "'1044.06"
$ws.Cells.Item(81, 8).Value = "'1942.34"
$ws.Cells.Item(82, 8).Value = "'6000.00"
$ws.Cells.Item(83, 8).Value = "'1960.00"
$ws.Cells.Item(84, 8).Value = "'1060.50"
$ws.Cells.Item(85, 8).Value = "'4470.00"
$ws.Cells.Item(86, 8).Value = "'1090.50"
$ws.Cells.Item(87, 8).Value = "'2127.00"
$ws.Cells.Item(88, 8).Value = "'177.15"
$ws.Cells.Item(89, 8).Value = "'509.80"
$ws.Cells.Item(90, 8).Value = "'29.50"
$ws.Cells.Item(91, 8).Value = "'142.27"
$ws.Cells.Item(92, 8).Value = "'23.79"
$ws.Cells.Item(93, 8).Value = "'4035.00"
$ws.Cells.Item(94, 8).Value = "'372.00"
$ws.Cells.Item(95, 8).Value = "'122.00"
$ws.Cells.Item(96, 8).Value = "'756.00"
$ws.Cells.Item(97, 8).Value = "'24561.15"
$ws.Cells.Item(98, 8).Value = "'12084.00"
$ws.Cells.Item(99, 8).Value = "'100.00"
$ws.Cells.Item(100, 8).Value = "'62.00"
$ws.Cells.Item(101, 8).Value = "'565.00"
$ws.Cells.Item(102, 8).Value = "'2882.00"
$ws.Cells.Item(103, 8).Value = "'33880.00"
$ws.Cells.Item(104, 8).Value = "'975.00"
$ws.Cells.Item(105, 8).Value = "'151.00"
$ws.Cells.Item(106, 8).Value = "'901.00"
$ws.Cells.Item(107, 8).Value = "'1902.49"
$ws.Cells.Item(108, 8).Value = "'4970.84"
$ws.Cells.Item(109, 8).Value = "'3630.00"
$ws.Cells.Item(110, 8).Value = "'250.00"
$ws.Cells.Item(111, 8).Value = "'500.00"
$ws.Cells.Item(112, 8).Value = "'1400.00"
$ws.Cells.Item(113, 8).Value = "'4772.24"
$ws.Cells.Item(114, 8).Value = "'290.00"
$ws.Cells.Item(115, 8).Value = "'500.00"
$ws.Cells.Item(116, 8).Value = "'150.00"
$ws.Cells.Item(117, 8).Value = "'1010.00"
$ws.Cells.Item(118, 8).Value = "'13050.04"
$ws.Cells.Item(119, 8).Value = "'350.00"
$ws.Cells.Item(120, 8).Value = "'2000.00"
$ws.Cells.Item(121, 8).Value = "'550.00"
$ws.Cells.Item(122, 8).Value = "'120.00"
$ws.Cells.Item(123, 8).Value = "'3250.00"
$ws.Cells.Item(124, 8).Value = "'7016.60"
$ws.Cells.Item(125, 8).Value = "'344.00"
$ws.Cells.Item(126, 8).Value = "'534.00"
$ws.Cells.Item(127, 8).Value = "'60.00"
$ws.Cells.Item(128, 8).Value = "'350.00"
$ws.Cells.Item(129, 8).Value = "'726.00"
$ws.Cells.Item(130, 8).Value = "'1450.00"
$ws.Cells.Item(131, 8).Value = "'1320.00"
$ws.Cells.Item(132, 8).Value = "'108900.00"
$ws.Cells.Item(133, 8).Value = "'10.50"
$ws.Cells.Item(134, 8).Value = "'950.00"
$ws.Cells.Item(135, 8).Value = "'4714.00"
$ws.Cells.Item(136, 8).Value = "'170.00"
$ws.Cells.Item(137, 8).Value = "'4706.00"
$ws.Cells.Item(138, 8).Value = "'1519.00"
$ws.Cells.Item(139, 8).Value = "'97.00"
$ws.Cells.Item(140, 8).Value = "'950.00"
$ws.Cells.Item(141, 8).Value = "'850.00"
$ws.Cells.Item(142, 8).Value = "'10108.19"
$ws.Cells.Item(143, 8).Value = "'1173.00"
$ws.Cells.Item(144, 8).Value = "'70.95"
$ws.Cells.Item(145, 8).Value = "'894.00"
$ws.Cells.Item(146, 8).Value = "'421.74"
$ws.Cells.Item(147, 8).Value = "'85.50"
$ws.Cells.Item(148, 8).Value = "'747.00"
$ws.Cells.Item(149, 8).Value = "'6441.00"
$ws.Cells.Item(150, 8).Value = "'1409.00"
$ws.Cells.Item(151, 8).Value = "'649.56"
$ws.Cells.Item(152, 8).Value = "'1750.00"
$ws.Cells.Item(153, 8).Value = "'50500.00"
$ws.Cells.Item(154, 8).Value = "'1766.00"
$ws.Cells.Item(155, 8).Value = "'356792.70"
$ws.Cells.Item(156, 8).Value = "'42000.00"
$ws.Cells.Item(157, 8).Value = "'908.26"
$ws.Cells.Item(158, 8).Value = "'6700.00"
$ws.Cells.Item(159, 8).Value = "'3670.00"
$ws.Cells.Item(160, 8).Value = "'2000.00"
$ws.Cells.Item(161, 8).Value = "'1542.95"
$ws.Cells.Item(162, 8).Value = "'360.00"
